# Inserts a new weekly record for "Rabanito" (Vega Central Mapocho de Santiago)
# at row 198, pushing the existing rows 198:236 down to 199:237.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 198 (shifts rows 198:236 -> 199:237)
$ws.Rows("198:198").Insert()

# Populate the newly inserted row with the new record's values
$ws.Range("A198").Value = 9
$ws.Range("B198").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C198").Value = "Metropolitana"
$ws.Range("D198").Value = 44637
$ws.Range("E198").Value = 13
$ws.Range("F198").Value = 300000001
$ws.Range("G198").Value = "Rabanito"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 5200
$ws.Range("K198").Value = 3000
$ws.Range("L198").Value = 3000
$ws.Range("M198").Value = 3000
$ws.Range("N198").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O198").Value = "Provincia de Chacabuco"
$ws.Range("P198").Value = 30
$ws.Range("Q198").Value = 100
$ws.Range("R198").Value = "Hortaliza"

# Apply the same date number format used by the rest of column D to the new cell
$ws.Range("D198").NumberFormat = $ws.Range("D199").NumberFormat
